# Added harvester and experiment design
# - Column B (harvester): "Retrofitted_0769" -> "S.GISH" for all data rows
# - Column D (experimentDesign): new value "90minuteInduction" for all data rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 22 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "S.GISH"
    $ws.Cells.Item($r, 4).Value = "90minuteInduction"
}
